$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.63
$summary.Range("B4").Value = -0.37
$summary.Range("B5").Value = -0.82
$summary.Range("B6").Value = 9
$summary.Range("B7").Value = 3
$summary.Range("B9").Value = 33.33

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.63
$status.Range("D4").Value = 9
$status.Range("E4").Value = -0.37
$status.Range("F4").Value = -0.37
$status.Range("G4").Value = 33.33

# --- New trade (#9) appended as row 10 on both "All Trades" and
#     "MarketMaking" sheets (they mirror each other) ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A10").Value = 9

    # "2026-02-17" looks like a date, so Excel's automatic type
    # detection would otherwise coerce it into a date serial number.
    # Force text formatting for the assignment, then restore the
    # cell's style to Normal so no stray formatting is left behind.
    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = "2026-02-17"
    $ws.Range("B10").Style = "Normal"

    $ws.Range("C10").Value = "13:08:22"
    $ws.Range("D10").Value = "MarketMaking"
    $ws.Range("E10").Value = "DOWN"
    $ws.Range("F10").Value = 0.8100000000000001
    $ws.Range("G10").Value = 0.9
    $ws.Range("H10").Value = "CLOSED"
    $ws.Range("I10").Value = 11.1111
    $ws.Range("J10").Value = 0.09
    $ws.Range("K10").Value = 99.63
    $ws.Range("L10").Value = 0
    $ws.Range("M10").Value = 0
    $ws.Range("N10").Value = 0.6
    $ws.Range("O10").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P10").Value = "early_exit"
    $ws.Range("Q10").Value = 0.13
}
